$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update FromAirport first (new unique string, gets appended after old "Tan Son Nhat..." is purged)
$ws.Range("E2").Value = "Da Nang International Airport"

# Update FlightNumber / AirplaneCode (swap-fix values)
$ws.Range("A2").Value = "QH3456"
$ws.Range("B2").Value = "QH1111"

# ToAirport stays "Noi Bai International Airport" - leave F2 untouched

# Departure time, duration, and prices
$ws.Range("C2").Value = 45602.333333333336
$ws.Range("D2").Value = 75
$ws.Range("G2").Value = 1990000
$ws.Range("H2").Value = 3400000
$ws.Range("I2").Value = 80000000

# Update the active selection to C3
$ws.Range("C3").Select() | Out-Null
